# Increment the "Förändrad" (Changed) date in column C by one day
# for every data row on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row from column C (date column), falling back to
# the worksheet's UsedRange if that lookup fails for some reason.
$xlUp = -4162
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End($xlUp).Row
if ($lastRow -lt 2) {
    $ur = $ws.UsedRange
    $lastRow = $ur.Row + $ur.Rows.Count - 1
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $v = $cell.Value2
    if ($v -ne $null) {
        $cell.Value2 = $v + 1
    }
}
